$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2526.4285
$ws.Range("J100").Value = 2914.5454
$ws.Range("L100").Value = 2914.5454
$ws.Range("N100").Value = -3996.5454

$ws.Range("H118").Value = 454.44446
$ws.Range("I118").Value = 278
$ws.Range("J118").Value = 675
$ws.Range("K118").Value = 834
$ws.Range("L118").Value = 2025
$ws.Range("M118").Value = 823
$ws.Range("N118").Value = -5339

$ws.Range("H138").Value = 3614.5764
$ws.Range("I138").Value = 2167.6667
$ws.Range("J138").Value = 4089.3438
$ws.Range("K138").Value = 6503.000100000001
$ws.Range("L138").Value = 12268.0314
$ws.Range("M138").Value = -1363.000100000001
$ws.Range("N138").Value = -22548.0314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4594.6665
$ws.Range("I3").Value = 7235
$ws.Range("J3").Value = 4066.6
$ws.Range("K3").Value = 7235
$ws.Range("L3").Value = 4066.6
$ws.Range("M3").Value = -7120
$ws.Range("N3").Value = -4296.6

$ws.Range("H32").Value = 7996.42
$ws.Range("I32").Value = 7649.408
$ws.Range("K32").Value = 7649.408
$ws.Range("M32").Value = -7362.408

$ws.Range("H97").Value = 999.65625
$ws.Range("I97").Value = 831.5599999999999
$ws.Range("J97").Value = 1600
$ws.Range("K97").Value = 831.5599999999999
$ws.Range("L97").Value = 1600
$ws.Range("M97").Value = -335.5599999999999
$ws.Range("N97").Value = -2592

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1678.1621
$ws.Range("I94").Value = 1572.4231
$ws.Range("J94").Value = 1928.091
$ws.Range("K94").Value = 1572.4231
$ws.Range("L94").Value = 1928.091
$ws.Range("M94").Value = -1121.4231
$ws.Range("N94").Value = -2830.091

$ws.Range("H99").Value = 2013.8823
$ws.Range("I99").Value = 1971.3334
$ws.Range("J99").Value = 2333
$ws.Range("K99").Value = 1971.3334
$ws.Range("L99").Value = 2333
$ws.Range("M99").Value = -473.3334
$ws.Range("N99").Value = -5329

$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344

$ws.Range("H107").Value = 2689.6667
$ws.Range("I107").Value = 2484.6667
$ws.Range("J107").Value = 3202.1667
$ws.Range("K107").Value = 2484.6667
$ws.Range("L107").Value = 3202.1667
$ws.Range("M107").Value = -564.6667000000002
$ws.Range("N107").Value = -7042.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2459988
$ws.Range("I58").Value = 3638606.2
$ws.Range("K58").Value = 3638606.2
$ws.Range("M58").Value = -3638403.2

$ws.Range("H110").Value = 41702
$ws.Range("J110").Value = 41702
$ws.Range("L110").Value = 41702
$ws.Range("N110").Value = -49882

$ws.Range("H136").Value = 2459988
$ws.Range("I136").Value = 3638606.2
$ws.Range("K136").Value = 10915818.6
$ws.Range("M136").Value = -10913268.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 949.9
$ws.Range("I4").Value = 333
$ws.Range("J4").Value = 1214.2858
$ws.Range("K4").Value = 999
$ws.Range("L4").Value = 3642.8574
$ws.Range("M4").Value = -887
$ws.Range("N4").Value = -3866.8574

$ws.Range("H17").Value = 230.66667
$ws.Range("J17").Value = 230.66667
$ws.Range("L17").Value = 692.00001
$ws.Range("N17").Value = -1030.00001

$ws.Range("H23").Value = 898.2222
$ws.Range("J23").Value = 334.66666
$ws.Range("L23").Value = 1003.99998
$ws.Range("N23").Value = -1473.99998

$ws.Range("H113").Value = 799.4516
$ws.Range("I113").Value = 802.6892
$ws.Range("J113").Value = 786.8421
$ws.Range("K113").Value = 2408.0676
$ws.Range("L113").Value = 2360.5263
$ws.Range("M113").Value = -238.0676000000003
$ws.Range("N113").Value = -6700.5263

$ws.Range("H134").Value = 4622.5
$ws.Range("I134").Value = 5205.75
$ws.Range("J134").Value = 4185.0625
$ws.Range("K134").Value = 15617.25
$ws.Range("L134").Value = 12555.1875
$ws.Range("M134").Value = -10547.25
$ws.Range("N134").Value = -22695.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 510000
$ws.Range("I10").Value = 510000
$ws.Range("K10").Value = 510000
$ws.Range("M10").Value = -509831

$ws.Range("H132").Value = 38790.434
$ws.Range("I132").Value = 78928.69500000001
$ws.Range("K132").Value = 236786.085
$ws.Range("M132").Value = -234256.085

$ws.Range("H140").Value = 53860
$ws.Range("J140").Value = 53860
$ws.Range("L140").Value = 53860
$ws.Range("N140").Value = -64220

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10000
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10224

$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10340

$ws.Range("H61").Value = 2014002.4
$ws.Range("J61").Value = 35468.332
$ws.Range("L61").Value = 35468.332
$ws.Range("N61").Value = -35872.332

$ws.Range("H68").Value = 5700.6665
$ws.Range("I68").Value = 5700.6665
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5700.6665
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4951.6665
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 5700.6665
$ws.Range("I71").Value = 5700.6665
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 28503.3325
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -24759.3325
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 2234.9285
$ws.Range("I82").Value = 1397.7778
$ws.Range("J82").Value = 3741.8
$ws.Range("K82").Value = 1397.7778
$ws.Range("L82").Value = 3741.8
$ws.Range("M82").Value = -1036.7778
$ws.Range("N82").Value = -4463.8

$ws.Range("H85").Value = 2234.9285
$ws.Range("I85").Value = 1397.7778
$ws.Range("J85").Value = 3741.8
$ws.Range("K85").Value = 1397.7778
$ws.Range("L85").Value = 3741.8
$ws.Range("M85").Value = -149.7778000000001
$ws.Range("N85").Value = -6237.8

$ws.Range("H113").Value = 2014002.4
$ws.Range("J113").Value = 35468.332
$ws.Range("L113").Value = 35468.332
$ws.Range("N113").Value = -39808.332

$ws.Range("H122").Value = 5824.9424
$ws.Range("I122").Value = 4740.2188
$ws.Range("J122").Value = 7560.5
$ws.Range("K122").Value = 14220.6564
$ws.Range("L122").Value = 22681.5
$ws.Range("M122").Value = -11770.6564
$ws.Range("N122").Value = -27581.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2993.3333
$ws.Range("I96").Value = 2993.3333
$ws.Range("K96").Value = 2993.3333
$ws.Range("M96").Value = -1620.3333

$ws.Range("H122").Value = 4158.4
$ws.Range("I122").Value = 2431.2144
$ws.Range("K122").Value = 7293.6432
$ws.Range("M122").Value = -4843.6432

$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
